$d = $word.ActiveDocument

# --- Edit 1: Trucha Fly Fishing Website bullet ---
# Replace the middle portion "Node.JS, Express.JS, Axios, and React.JS"
# with "HTML, CSS, and Javascript" inside the existing sentence. This causes
# Word to split the original run into three runs (before / new text / after)
# while preserving the original run formatting.
$d.Content.Find.Execute(
    "Node.JS, Express.JS, Axios, and React.JS",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "HTML, CSS, and Javascript", 2)

# --- Edit 2: MySQL Database bullet ---
# Replace the whole previous sentence with two sentences describing
# current work, split so the formatting-preserving run split matches
# what Word naturally produces.
$d.Content.Find.Execute(
    "Utilized MySQL Database to create a functional blog posting feature so owner could display his trips",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Currently working with MySQL, Node.JS, Express.JS and React.JS to implement a blog posting feature", 2)
